$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $r = $t.Rows.Item($i)
    if ($r.Cells.Item(1).Range.Text -like "Overall Climate Vulnerability*") {
        $r.Delete()
        break
    }
}
